{"js": "// Update the delivery date from \"03/07\" to \"11/09\" and remove the\n// stray \"_GoBack\" bookmark left over from the previous save.\n\n// 1) Update the date text \"03/07\" -> \"11/09\".\nconst body = context.document.body;\nconst searchResults = body.search(\"03/07\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"11/09\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Remove the \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd pair).\nconst goBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\ngoBack.load(\"isNullObject\");\nawait context.sync();\n\nif (!goBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Update the delivery date from \"03/07\" to \"11/09\" and remove the\n# stray \"_GoBack\" bookmark left over from the previous save.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date text \"03/07\" -> \"11/09\".\n$find = $d.Content.Find\n$find.Text = \"03/07\"\n$find.Replacement.Text = \"11/09\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 2) Remove the \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd pair).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
